# WAT/RRC project refactor: trim the "Test Suite" sheet down to just the
# header row plus the WAT and WATIAM module rows, dropping every other
# module/suite row (IAM, Search, Watchlist, Notification, Authoring,
# Profile, Sanity, ENW, LocalRun, DRA, IPA, RCC, customercare, PUBLONS, ...).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Test Suite")

# Row 19 (PUBLONS) is below the rows we want to keep - remove it first so
# the subsequent bulk delete doesn't need to account for the shift.
$ws.Rows("19:19").Delete()

# Rows 2-16 hold every module row except WAT (old row 17) and WATIAM
# (old row 18). Deleting them shifts WAT/WATIAM up to become the new
# rows 2 and 3, directly under the header row.
$ws.Rows("2:16").Delete()

# Match the saved selection state from the edited workbook.
$ws.Range("C3").Select()
